# Adding the changes we made on may 9th
#
# The sheet holds an (x, y, z) gyroscope sample table in A1:C21 (A1:C1 are
# the headers "x","y","z"). This edit:
#   1. Drops the first three data rows (old rows 2-4), so the old row 5
#      becomes the new row 2, etc.
#   2. Appends 13 freshly captured samples after the (now shifted) data,
#      extending the table down to row 31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the first three samples (old rows 2-4); rows below shift up.
$ws.Range("A2:C4").EntireRow.Delete()

# 2. Full data block for rows 2-31 after the shift + append, written in one
#    shot so there is no ambiguity about which old rows survived.
$data = @(
    @(-0.5190548883925257, -1.450479206421397, -0.2456109922566745),
    @(-0.2103044650108719, -1.79173934116638, -0.6462500674261389),
    @(2.22406245532224, -1.244166610052261, -0.3580538825594288),
    @(1.525593809515442, -1.465595992349037, 0.216453205767295),
    @(-4.037772781986149, -3.096113548004377, 2.077658006398807),
    @(-7.127989895909813, 1.020252332412979, 0.7169627440061292),
    @(0.1363662918695203, -0.7235090681117633, 1.348089827050427),
    @(6.111373266727803, -3.888724038926818, -0.6717052065211169),
    @(2.017424072293039, -1.605387079201143, -3.74285111324393),
    @(-3.043187412426637, -0.6237803290216224, 1.956819472553102),
    @(-10.0801253147263, -3.291420891988214, 5.852694562870816),
    @(-1.435964323633735, -3.116479964565011, -1.45930368780239),
    @(4.520515643435422, -6.424601422797, -1.258769486447896),
    @(3.804642815384075, -2.879516224209357, 0.2698007882070206),
    @(2.50941535251603, 3.106572741227184, 1.131389019300618),
    @(-5.14737601803366, 0.1312354588679825, 5.117585698477669),
    @(-2.502922200470437, 3.250862176469831, -2.965100693188241),
    @(2.697737143194075, -2.115886496125361, -4.09576786507806),
    @(5.033141990359741, -4.757362827122633, -1.592488460403557),
    @(4.163411040957872, 2.378721077665152, -1.909676737064927),
    @(-2.440776584817322, 2.270052646561497, 2.502104572254968),
    @(-4.33349856712834, 5.117880192592036, -1.15704054626636),
    @(0.3903590963899126, 12.85140922601278, -6.864528203182092),
    @(2.438172858395985, 7.542828546153372, -8.415844210617843),
    @(-0.09949780518202356, -6.494732479397333, -8.307635595472599),
    @(-1.967174616863385, -6.35848597656913, -2.059108837045416),
    @(-9.929815484465564, -8.435556281384915, 4.079216970814213),
    @(-5.980313006064948, -16.70210841062257, 0.5289435215133755),
    @(5.703074371214343, -13.00758359758107, -1.666179502181835),
    @(3.259694958762311, -1.829939403122192, -3.400759592759585)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
